$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.054997697132236
$ws.Cells.Item(2, 4).Value = 1.054771049449333
$ws.Cells.Item(2, 5).Value = 1.068480997652134
$ws.Cells.Item(2, 6).Value = 1.075517826583241
$ws.Cells.Item(2, 9).Value = 1.045925623614579
$ws.Cells.Item(2, 10).Value = 1.060006862798617
$ws.Cells.Item(2, 11).Value = 1.057513342850015
$ws.Cells.Item(2, 12).Value = 1.071186065228156
$ws.Cells.Item(2, 13).Value = 1.07820418588485
$ws.Cells.Item(2, 14).Value = 1.023669340828943

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.056201279661739
$ws.Cells.Item(3, 4).Value = 1.055702721894697
$ws.Cells.Item(3, 5).Value = 1.069709257235593
$ws.Cells.Item(3, 6).Value = 1.076888143554927
$ws.Cells.Item(3, 9).Value = 1.046275464233516
$ws.Cells.Item(3, 10).Value = 1.060860719938657
$ws.Cells.Item(3, 11).Value = 1.058258022143918
$ws.Cells.Item(3, 12).Value = 1.072229287025454
$ws.Cells.Item(3, 13).Value = 1.079390478533576
$ws.Cells.Item(3, 14).Value = 1.023963296652942

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.05697970682106
$ws.Cells.Item(4, 4).Value = 1.056305166782596
$ws.Cells.Item(4, 5).Value = 1.070504041051371
$ws.Cells.Item(4, 6).Value = 1.077775121265602
$ws.Cells.Item(4, 9).Value = 1.046500435295583
$ws.Cells.Item(4, 10).Value = 1.061412311453375
$ws.Cells.Item(4, 11).Value = 1.058738843530167
$ws.Cells.Item(4, 12).Value = 1.0729037699741
$ws.Cells.Item(4, 13).Value = 1.080157837369014
$ws.Cells.Item(4, 14).Value = 1.02415297318758

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.057306870971008
$ws.Cells.Item(5, 4).Value = 1.056558337616408
$ws.Cells.Item(5, 5).Value = 1.070838174367047
$ws.Cells.Item(5, 6).Value = 1.078148078756744
$ws.Cells.Item(5, 9).Value = 1.046594678768652
$ws.Cells.Item(5, 10).Value = 1.061643983806342
$ws.Cells.Item(5, 11).Value = 1.058940733951567
$ws.Cells.Item(5, 12).Value = 1.073187192373587
$ws.Cells.Item(5, 13).Value = 1.080480376805211
$ws.Cells.Item(5, 14).Value = 1.024232586049705

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.057361798256651
$ws.Cells.Item(6, 4).Value = 1.056600840455455
$ws.Cells.Item(6, 5).Value = 1.070894277241325
$ws.Cells.Item(6, 6).Value = 1.078210704255952
$ws.Cells.Item(6, 9).Value = 1.046610483073156
$ws.Cells.Item(6, 10).Value = 1.061682869946253
$ws.Cells.Item(6, 11).Value = 1.058974617796196
$ws.Cells.Item(6, 12).Value = 1.07323477265564
$ws.Cells.Item(6, 13).Value = 1.080534529242034
$ws.Cells.Item(6, 14).Value = 1.024245945953053

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.056984078742024
$ws.Cells.Item(7, 4).Value = 1.056308550043512
$ws.Cells.Item(7, 5).Value = 1.070508505730283
$ws.Cells.Item(7, 6).Value = 1.077780104455699
$ws.Cells.Item(7, 9).Value = 1.046501695894529
$ws.Cells.Item(7, 10).Value = 1.061415407920493
$ws.Cells.Item(7, 11).Value = 1.058741542168582
$ws.Cells.Item(7, 12).Value = 1.072907557587606
$ws.Cells.Item(7, 13).Value = 1.080162147382915
$ws.Cells.Item(7, 14).Value = 1.024154037477743

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.055404531696952
$ws.Cells.Item(8, 4).Value = 1.055085997924419
$ws.Cells.Item(8, 5).Value = 1.068896091234048
$ws.Cells.Item(8, 6).Value = 1.075980872808761
$ws.Cells.Item(8, 9).Value = 1.046044143974817
$ws.Cells.Item(8, 10).Value = 1.060295617327058
$ws.Cells.Item(8, 11).Value = 1.057765225865576
$ws.Cells.Item(8, 12).Value = 1.071538742127819
$ws.Cells.Item(8, 13).Value = 1.078605152430765
$ws.Cells.Item(8, 14).Value = 1.02376879504594

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.05261822424348
$ws.Cells.Item(9, 4).Value = 1.052928521903111
$ws.Cells.Item(9, 5).Value = 1.066054858881791
$ws.Cells.Item(9, 6).Value = 1.072812517975681
$ws.Cells.Item(9, 9).Value = 1.045227135102106
$ws.Cells.Item(9, 10).Value = 1.058315361365682
$ws.Cells.Item(9, 11).Value = 1.056036847863565
$ws.Cells.Item(9, 12).Value = 1.069122413188943
$ws.Cells.Item(9, 13).Value = 1.075859506642642
$ws.Cells.Item(9, 14).Value = 1.023085854678582

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.050758563856461
$ws.Cells.Item(10, 4).Value = 1.051487996463275
$ws.Cells.Item(10, 5).Value = 1.064160611108945
$ws.Cells.Item(10, 6).Value = 1.070701555049879
$ws.Cells.Item(10, 9).Value = 1.044675199042721
$ws.Cells.Item(10, 10).Value = 1.056990365296654
$ws.Cells.Item(10, 11).Value = 1.054879155850927
$ws.Cells.Item(10, 12).Value = 1.067508522968162
$ws.Cells.Item(10, 13).Value = 1.07402758457129
$ws.Cells.Item(10, 14).Value = 1.022627786055419

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.049952774579146
$ws.Cells.Item(11, 4).Value = 1.050863693376134
$ws.Cells.Item(11, 5).Value = 1.0633403302949
$ws.Cells.Item(11, 6).Value = 1.069787752233937
$ws.Cells.Item(11, 9).Value = 1.044434472567105
$ws.Cells.Item(11, 10).Value = 1.05641546197446
$ws.Cells.Item(11, 11).Value = 1.054376555448147
$ws.Cells.Item(11, 12).Value = 1.066808952422551
$ws.Cells.Item(11, 13).Value = 1.073233958198101
$ws.Cells.Item(11, 14).Value = 1.022428773760097

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.049653383934834
$ws.Cells.Item(12, 4).Value = 1.050631715901339
$ws.Cells.Item(12, 5).Value = 1.063035630052221
$ws.Cells.Item(12, 6).Value = 1.069448360875761
$ws.Cells.Item(12, 9).Value = 1.044344794581942
$ws.Cells.Item(12, 10).Value = 1.056201739423616
$ws.Cells.Item(12, 11).Value = 1.054189668622721
$ws.Cells.Item(12, 12).Value = 1.066548986474814
$ws.Cells.Item(12, 13).Value = 1.072939108441916
$ws.Cells.Item(12, 14).Value = 1.022354751239367

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.04971760812048
$ws.Cells.Item(13, 4).Value = 1.050681479659906
$ws.Cells.Item(13, 5).Value = 1.063100989846772
$ws.Cells.Item(13, 6).Value = 1.069521159908404
$ws.Cells.Item(13, 9).Value = 1.044364042659376
$ws.Cells.Item(13, 10).Value = 1.056247591715162
$ws.Cells.Item(13, 11).Value = 1.054229765495417
$ws.Cells.Item(13, 12).Value = 1.066604755262446
$ws.Cells.Item(13, 13).Value = 1.073002357545591
$ws.Cells.Item(13, 14).Value = 1.022370633876939

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.049928028598239
$ws.Cells.Item(14, 4).Value = 1.050844519757836
$ws.Cells.Item(14, 5).Value = 1.063315143925676
$ws.Cells.Item(14, 6).Value = 1.069759697328052
$ws.Cells.Item(14, 9).Value = 1.044427065095445
$ws.Cells.Item(14, 10).Value = 1.056397799240952
$ws.Cells.Item(14, 11).Value = 1.054361111386257
$ws.Cells.Item(14, 12).Value = 1.066787465901235
$ws.Cells.Item(14, 13).Value = 1.073209587097158
$ws.Cells.Item(14, 14).Value = 1.022422657088378

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.050057664298597
$ws.Cells.Item(15, 4).Value = 1.050944963033011
$ws.Cells.Item(15, 5).Value = 1.063447089731081
$ws.Cells.Item(15, 6).Value = 1.069906672742777
$ws.Cells.Item(15, 9).Value = 1.044465860627334
$ws.Cells.Item(15, 10).Value = 1.056490323430398
$ws.Cells.Item(15, 11).Value = 1.054442011554766
$ws.Cells.Item(15, 12).Value = 1.066900024741214
$ws.Cells.Item(15, 13).Value = 1.073337259822204
$ws.Cells.Item(15, 14).Value = 1.022454696974579

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.05081202972778
$ws.Cells.Item(16, 4).Value = 1.05152941777272
$ws.Cells.Item(16, 5).Value = 1.064215048983311
$ws.Cells.Item(16, 6).Value = 1.07076220630792
$ws.Cells.Item(16, 9).Value = 1.044691138648249
$ws.Cells.Item(16, 10).Value = 1.057028494898224
$ws.Cells.Item(16, 11).Value = 1.054912484009359
$ws.Cells.Item(16, 12).Value = 1.067554935259244
$ws.Cells.Item(16, 13).Value = 1.074080246451161
$ws.Cells.Item(16, 14).Value = 1.022640979780814

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.0512850755329
$ws.Cells.Item(17, 4).Value = 1.051895883555964
$ws.Cells.Item(17, 5).Value = 1.064696751926121
$ws.Cells.Item(17, 6).Value = 1.071298926831037
$ws.Cells.Item(17, 9).Value = 1.044831984610866
$ws.Cells.Item(17, 10).Value = 1.05736576072534
$ws.Cells.Item(17, 11).Value = 1.055207246582144
$ws.Cells.Item(17, 12).Value = 1.067965542125342
$ws.Cells.Item(17, 13).Value = 1.074546195481732
$ws.Cells.Item(17, 14).Value = 1.022757651522033

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.051560942850761
$ws.Cells.Item(18, 4).Value = 1.052109584260638
$ws.Cells.Item(18, 5).Value = 1.064977715581809
$ws.Cells.Item(18, 6).Value = 1.071612011940132
$ws.Cells.Item(18, 9).Value = 1.044913970380783
$ws.Cells.Item(18, 10).Value = 1.057562369293372
$ws.Cells.Item(18, 11).Value = 1.055379050107413
$ws.Cells.Item(18, 12).Value = 1.068204970397423
$ws.Cells.Item(18, 13).Value = 1.074817938026074
$ws.Cells.Item(18, 14).Value = 1.022825639992358

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.051654997693503
$ws.Cells.Item(19, 4).Value = 1.052182441802771
$ws.Cells.Item(19, 5).Value = 1.065073516067939
$ws.Cells.Item(19, 6).Value = 1.071718770222778
$ws.Cells.Item(19, 9).Value = 1.044941897058952
$ws.Cells.Item(19, 10).Value = 1.057629388644
$ws.Cells.Item(19, 11).Value = 1.055437609259462
$ws.Cells.Item(19, 12).Value = 1.068286597147303
$ws.Cells.Item(19, 13).Value = 1.074910588907238
$ws.Cells.Item(19, 14).Value = 1.022848811432252

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.05123432762632
$ws.Cells.Item(20, 4).Value = 1.051856570696434
$ws.Cells.Item(20, 5).Value = 1.064645070349898
$ws.Cells.Item(20, 6).Value = 1.071241339224846
$ws.Cells.Item(20, 9).Value = 1.044816890483065
$ws.Cells.Item(20, 10).Value = 1.057329586984513
$ws.Cells.Item(20, 11).Value = 1.055175634443659
$ws.Cells.Item(20, 12).Value = 1.067921495335397
$ws.Cells.Item(20, 13).Value = 1.074496207474556
$ws.Cells.Item(20, 14).Value = 1.022745140390997

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.049866067371557
$ws.Cells.Item(21, 4).Value = 1.050796510829831
$ws.Cells.Item(21, 5).Value = 1.06325208122306
$ws.Cells.Item(21, 6).Value = 1.069689453042752
$ws.Cells.Item(21, 9).Value = 1.044408513783152
$ws.Cells.Item(21, 10).Value = 1.056353571793414
$ws.Cells.Item(21, 11).Value = 1.054322438791321
$ws.Cells.Item(21, 12).Value = 1.066733665345817
$ws.Cells.Item(21, 13).Value = 1.073148564880949
$ws.Cells.Item(21, 14).Value = 1.022407340327393

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.049005296769669
$ws.Cells.Item(22, 4).Value = 1.050129525214581
$ws.Cells.Item(22, 5).Value = 1.0623761862962
$ws.Cells.Item(22, 6).Value = 1.068713925595589
$ws.Cells.Item(22, 9).Value = 1.044150238024279
$ws.Cells.Item(22, 10).Value = 1.055738882326033
$ws.Cells.Item(22, 11).Value = 1.053784850690353
$ws.Cells.Item(22, 12).Value = 1.065986167042483
$ws.Cells.Item(22, 13).Value = 1.072300890847008
$ws.Cells.Item(22, 14).Value = 1.022194370347384

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.049461655067631
$ws.Cells.Item(23, 4).Value = 1.050483153270827
$ws.Cells.Item(23, 5).Value = 1.062840522011997
$ws.Cells.Item(23, 6).Value = 1.069231052751925
$ws.Cells.Item(23, 9).Value = 1.044287298605487
$ws.Cells.Item(23, 10).Value = 1.056064839105281
$ws.Cells.Item(23, 11).Value = 1.054069945844621
$ws.Cells.Item(23, 12).Value = 1.066382493600457
$ws.Cells.Item(23, 13).Value = 1.072750293691368
$ws.Cells.Item(23, 14).Value = 1.022307325067157

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.051257258582446
$ws.Cells.Item(24, 4).Value = 1.051874334648008
$ws.Cells.Item(24, 5).Value = 1.064668423046655
$ws.Cells.Item(24, 6).Value = 1.071267360507651
$ws.Cells.Item(24, 9).Value = 1.044823711386635
$ws.Cells.Item(24, 10).Value = 1.057345932689712
$ws.Cells.Item(24, 11).Value = 1.055189918999193
$ws.Cells.Item(24, 12).Value = 1.067941398405549
$ws.Cells.Item(24, 13).Value = 1.074518795019585
$ws.Cells.Item(24, 14).Value = 1.022750793830917

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.053338915047045
$ws.Cells.Item(25, 4).Value = 1.053486665349069
$ws.Cells.Item(25, 5).Value = 1.066789392723426
$ws.Cells.Item(25, 6).Value = 1.073631377956531
$ws.Cells.Item(25, 9).Value = 1.045439628787756
$ws.Cells.Item(25, 10).Value = 1.058828148833411
$ws.Cells.Item(25, 11).Value = 1.056484628544889
$ws.Cells.Item(25, 12).Value = 1.069747613404274
$ws.Cells.Item(25, 13).Value = 1.076569575074501
$ws.Cells.Item(25, 14).Value = 1.023262898596245
